$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'265.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.59%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'26.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.81%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'4.702"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.07%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.06081"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.82%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'6.741"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.37%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.8506"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.08%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.9058"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.88%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1412"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.11%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.04833"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'6.18%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07094"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.15%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.03180"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.45%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.09023"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.25%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.001536"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.24%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.0006071"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.39%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.005970"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.49%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.456"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.05%"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'3.172"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.27%"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'2.277"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'3.78%"
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'-1.65%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'0.1300"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-0.78%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'4.089"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.15%"
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'0.08%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001185"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-2.63%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004132"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'8.67%"
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'0.04%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.0001682"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'5.06%"
$ws.Range("E27").Style = "Normal"

$ws.Range("E40").Value = "'-0.59%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.1113"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.03%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.004169"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.08%"
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'-3.32%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.01261"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-8.68%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005114"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.43%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.05%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D48").Value = "'0.1367"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-18.50%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0002002"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.05%"
$ws.Range("E50").Style = "Normal"
